# Insert a new weekly price-report row before row 52 (Ají / Cristal,
# Región del Maule) shifting every subsequent row down by one, so the
# sheet grows from A1:R136 to A1:R137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push old row 52 (and everything below it) down one row.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44967
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112021
$ws.Range("G52").Value = "Ají"
$ws.Range("H52").Value = "Cristal"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 30
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = 15000
$ws.Range("N52").Value = "`$/saco 25 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 600
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
